# -----------------------------------------------------------------------
# "Update countries & provincias Spain" - refresh the COVID country table:
#   * Update the "Datos actualizados" timestamp in A1
#   * Refresh the per-country statistics (cols B:H) with new totals
#   * A handful of countries swapped rank (their totals crossed over),
#     so the country names in column A for those rows are corrected too
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 29 de Junio de 2020 a las 19:34'

$ws.Range("B4").Value = 2652334
$ws.Range("C4").Value = 15257
$ws.Range("D4").Value = 1099188
$ws.Range("E4").Value = 1424589
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = 128557

$ws.Range("B5").Value = 1352708
$ws.Range("C5").Value = 7454
$ws.Range("E5").Value = 561086
$ws.Range("G5").Value = 116
$ws.Range("H5").Value = 57774

$ws.Range("B7").Value = 566931
$ws.Range("C7").Value = 17734
$ws.Range("D7").Value = 334907
$ws.Range("E7").Value = 215125
$ws.Range("G7").Value = 412
$ws.Range("H7").Value = 16899

$ws.Range("B9").Value = 296050
$ws.Range("C9").Value = 200
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 28346

$ws.Range("B11").Value = 275999
$ws.Range("C11").Value = 4017
$ws.Range("D11").Value = 236154
$ws.Range("E11").Value = 34270
$ws.Range("G11").Value = 66
$ws.Range("H11").Value = 5575

$ws.Range("B16").Value = 198613
$ws.Range("C16").Value = 1374
$ws.Range("D16").Value = 171809
$ws.Range("E16").Value = 21689
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 5115

$ws.Range("B17").Value = 195206
$ws.Range("C17").Value = 342
$ws.Range("E17").Value = 8074
$ws.Range("G17").Value = 3
$ws.Range("H17").Value = 9032

$ws.Range("B31").Value = 55665
$ws.Range("C31").Value = 410
$ws.Range("D31").Value = 27430
$ws.Range("E31").Value = 23733
$ws.Range("G31").Value = 73
$ws.Range("H31").Value = 4502

$ws.Range("B50").Value = 25462
$ws.Range("C50").Value = 23
$ws.Range("E50").Value = 363

$ws.Range("B53").Value = 24276
$ws.Range("C53").Value = 521
$ws.Range("D53").Value = 17189
$ws.Range("E53").Value = 6768

$ws.Range("E54").Value = 8131
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 188

$ws.Range("A59").Value = 'Azerbaiyan'
$ws.Range("B59").Value = 16968
$ws.Range("C59").Value = 544
$ws.Range("D59").Value = 9369
$ws.Range("E59").Value = 7393
$ws.Range("G59").Value = 8
$ws.Range("H59").Value = 206

$ws.Range("A60").Value = 'Guatemala'
$ws.Range("B60").Value = 16930
$ws.Range("C60").Value = 533
$ws.Range("D60").Value = 3152
$ws.Range("E60").Value = 13051
$ws.Range("G60").Value = 21
$ws.Range("H60").Value = 727

$ws.Range("B68").Value = 12290
$ws.Range("C68").Value = 238
$ws.Range("D68").Value = 8833
$ws.Range("E68").Value = 3232
$ws.Range("G68").Value = 4
$ws.Range("H68").Value = 225

$ws.Range("B95").Value = 3774
$ws.Range("C95").Value = 313
$ws.Range("D95").Value = 1352
$ws.Range("E95").Value = 2407
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = 15

$ws.Range("B105").Value = 2337
$ws.Range("C105").Value = 13
$ws.Range("D105").Value = 1927
$ws.Range("E105").Value = 402

$ws.Range("B123").Value = 1450
$ws.Range("C123").Value = 23
$ws.Range("D123").Value = 961
$ws.Range("E123").Value = 429

$ws.Range("B126").Value = 1172
$ws.Range("C126").Value = 3
$ws.Range("E126").Value = 93

$ws.Range("A129").Value = 'Yemen'
$ws.Range("C129").Value = 10
$ws.Range("D129").Value = 432
$ws.Range("E129").Value = 392
$ws.Range("G129").Value = 2
$ws.Range("H129").Value = 304

$ws.Range("A130").Value = 'Jordania'
$ws.Range("B130").Value = 1128
$ws.Range("C130").Value = 7
$ws.Range("D130").Value = 867
$ws.Range("E130").Value = 252
$ws.Range("H130").Value = 9

$ws.Range("B153").Value = 522
$ws.Range("C153").Value = 1
$ws.Range("E153").Value = 48

$ws.Range("A185").Value = 'Seychelles'
$ws.Range("C185").Value = 7
$ws.Range("D185").Value = 11
$ws.Range("E185").Value = 66
$ws.Range("H185").Value = 0

$ws.Range("A186").Value = 'San Martin (Parte Holandesa)'
$ws.Range("B186").Value = 77
$ws.Range("D186").Value = 62
$ws.Range("E186").Value = 0
$ws.Range("H186").Value = 15

$ws.Range("A187").Value = 'Butan'
$ws.Range("B187").Value = 76
$ws.Range("D187").Value = 44
$ws.Range("E187").Value = 32
